# Electricity Technology Shareweights.xlsx - EU update
# "updated CCS vars and ETS"
$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$etsSheet   = $wb.Worksheets.Item("ETS")

# --- ETS sheet: hydro (row 5) shareweights changed from 0 to 1 for all years (B5:AF5) ---
$etsSheet.Range("B5:AF5").Value = 1

# --- About sheet: merge the two "hydro/crude oil/fuel oil" footnote rows (A20:A21) into one ---
$aboutSheet.Range("A20").Value = "EU values are set to zero for new crude oil, fuel oil. There are some small scale hydro plants under construction, so this value was set to 1. Most (18 of 27) EU countries have committed to phasing out coal and could consider setting coal value to zero when all countries have committed to phase out all coal. "
[void]$aboutSheet.Range("A21").EntireRow.Delete()

# --- Update selections to match the saved state in the authored file ---
[void]$etsSheet.Range("B12").Select()
[void]$aboutSheet.Range("A21").Select()
[void]$aboutSheet.Activate()
